$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet: bump the generation Date and the concept Count
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-08-13T14:50:19+00:00"
$meta.Range("B22").Value = "8"

# ---------------------------------------------------------------------------
# Concepts sheet: add "Case-only" / "Case-parent-duo" right after the header,
# and append "Nuclear-family" at the end of the table.
# ---------------------------------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

# The table currently spans rows 1 (header) .. 6 (5 data rows). It needs to
# grow to rows 1 .. 9 (8 data rows). Give the three brand-new trailing rows
# (7, 8 and 9) the same formatting as the existing data rows before anything
# else touches them, by copying the format of the last existing data row.
$concepts.Range("A6:D6").Copy()
$concepts.Range("A7:D9").PasteSpecial(-4122)

# Shift the five existing data rows (2-6) down two rows (to 4-8), working
# from the bottom up so that source rows are not clobbered before they are
# read.
for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 2
    $concepts.Range("A$dest").Value = $concepts.Range("A$r").Value2
    $concepts.Range("B$dest").Value = $concepts.Range("B$r").Value2
    $concepts.Range("C$dest").Value = $concepts.Range("C$r").Value2
}

# Write the two newly introduced rows into the freed-up slots.
$concepts.Range("A2").Value = "1"
$concepts.Range("B2").Value = "Case-only"
$concepts.Range("C2").Value = "Case only"

$concepts.Range("A3").Value = "1"
$concepts.Range("B3").Value = "Case-parent-duo"
$concepts.Range("C3").Value = "Case-parent duo"

# Append the new trailing row (row 9 - already formatted above).
$concepts.Range("A9").Value = "1"
$concepts.Range("B9").Value = "Nuclear-family"
$concepts.Range("C9").Value = "Nuclear family"
